$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove the old "units" row (row 2) and the stray unlabeled data row
# (former row 5, which had no canton name) so the data rows are contiguous.
$ws.Rows.Item(5).Delete()
$ws.Rows.Item(2).Delete()

# Reset the leftover formatting on the cells that used to carry the
# "(m3/s)/(MW)/(GWh)" units-row styling so the new header row starts clean.
$ws.Range("A1:K1").Style = "Normal"

# Write the new single header row.
$ws.Range("A1").Value = "idx"
$ws.Range("B1").Value = "idx2"
$ws.Range("C1").Value = "Name"
$ws.Range("D1").Value = "Date Start"
$ws.Range("E1").Value = "Date End"
$ws.Range("F1").Value = "(m3/s)"
$ws.Range("G1").Value = "(MW1)"
$ws.Range("H1").Value = "(MW2)"
$ws.Range("I1").Value = "(GWh) Winter"
$ws.Range("J1").Value = "(GWh) Summer"
$ws.Range("K1").Value = "(GWh) Year"

# The unit-bearing columns (F:K) get a dedicated Arial 9 style distinct
# from the plain header cells (A:E stay on the default "Normal" style).
# Define it as a transient named style so the resulting cell format (xfId=0,
# fontId=1, applyFont only) is picked up, then drop the named style itself -
# the cells keep referencing the new cell-format record.
$unitStyle = $wb.Styles.Add("HeaderUnit")
$unitStyle.Font.Name = "Arial"
$unitStyle.Font.Size = 9
$ws.Range("F1:K1").Style = "HeaderUnit"
$wb.Styles.Item("HeaderUnit").Delete()

$ws.Range("A4:K4").Select()
